{"js": "// Week 41-2 log entry: add Victor's line after Andreas's line in the\n// \"Tirsdag\" (Tuesday) section. Mirrors a user putting the cursor at the\n// end of Andreas's paragraph, pressing Enter, and typing the new text\n// (same run formatting: sz=22 / szCs=22, i.e. 11pt).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// The paragraph that currently ends the document body is Andreas's line.\nconst andreasParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the new paragraph right after it, with Victor's text.\nconst victorParagraph = andreasParagraph.insertParagraph(\n  \"Victor: Jeg laver udgangspunkt til vores poster.\",\n  \"After\"\n);\n\n// Give the (now split) paragraph mark of Andreas's paragraph the same\n// 11pt (sz/szCs=22) formatting the run already carries - this is what\n// Word records when you press Enter at the end of an 11pt run.\nandreasParagraph.font.size = 11;\nandreasParagraph.font.sizeBidirectional = 11;\n\nawait context.sync();\n", "ps1": "# Week 41-2 log entry: add Victor's line after Andreas's line in the\n# \"Tirsdag\" (Tuesday) section. Mirrors a user putting the cursor at the\n# end of Andreas's paragraph, pressing Enter, and typing the new text\n# (same run formatting: sz=22 / szCs=22, i.e. 11pt).\n\n$d = $word.ActiveDocument\n\n# Locate Andreas's paragraph by its text.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"Andreas: Jeg laver problemformulering samt undersp\u00f8rgsm\u00e5l til vores poster\")\n\nif ($found) {\n    $andreasPara = $findRange.Paragraphs(1)\n\n    # Insert a new, empty paragraph right after Andreas's paragraph.\n    $insertionPoint = $andreasPara.Range\n    $insertionPoint.Collapse(0)  # wdCollapseEnd\n    $insertionPoint.InsertParagraphAfter()\n\n    # Fill the new paragraph with Victor's text.\n    $victorPara = $andreasPara.Next()\n    $victorRange = $victorPara.Range\n    $victorRange.Collapse(0)\n    $victorRange.Text = \"Victor: Jeg laver udgangspunkt til vores poster.\"\n\n    # Give the (now split) paragraph mark of Andreas's paragraph the same\n    # 11pt (sz/szCs=22) formatting the run already carries - this is what\n    # Word records when you press Enter at the end of an 11pt run.\n    $andreasPara.Range.Font.Size = 11\n    $andreasPara.Range.Font.SizeBi = 11\n}\n"}
